$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 39222216
$ws.Range("I132").Value = 31919854
$ws.Range("J132").Value = 125025000
$ws.Range("K132").Value = 95759562
$ws.Range("L132").Value = 375075000
$ws.Range("M132").Value = -95757032
$ws.Range("N132").Value = -375080060

$ws.Range("H137").Value = 1239
$ws.Range("I137").Value = 1130.7778
$ws.Range("J137").Value = 1543.375
$ws.Range("K137").Value = 3392.3334
$ws.Range("L137").Value = 4630.125
$ws.Range("M137").Value = -842.3334000000004
$ws.Range("N137").Value = -9730.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 886.52
$ws.Range("I61").Value = 813.5454999999999
$ws.Range("J61").Value = 1028.1765
$ws.Range("K61").Value = 813.5454999999999
$ws.Range("L61").Value = 1028.1765
$ws.Range("M61").Value = -601.5454999999999
$ws.Range("N61").Value = -1452.1765

$ws.Range("H74").Value = 640.9423
$ws.Range("I74").Value = 630.7826
$ws.Range("J74").Value = 718.8333
$ws.Range("K74").Value = 630.7826
$ws.Range("L74").Value = 718.8333
$ws.Range("M74").Value = 243.2174
$ws.Range("N74").Value = -2466.8333

$ws.Range("H77").Value = 640.9423
$ws.Range("I77").Value = 630.7826
$ws.Range("J77").Value = 718.8333
$ws.Range("K77").Value = 3153.913
$ws.Range("L77").Value = 3594.1665
$ws.Range("M77").Value = 1214.087
$ws.Range("N77").Value = -12330.1665

$ws.Range("H132").Value = 16044186
$ws.Range("I132").Value = 20409102
$ws.Range("J132").Value = 3462957.2
$ws.Range("K132").Value = 61227306
$ws.Range("L132").Value = 10388871.6
$ws.Range("M132").Value = -61224776
$ws.Range("N132").Value = -10393931.6

$ws.Range("H136").Value = 886.52
$ws.Range("I136").Value = 813.5454999999999
$ws.Range("J136").Value = 1028.1765
$ws.Range("K136").Value = 2440.6365
$ws.Range("L136").Value = 3084.5295
$ws.Range("M136").Value = 109.3635000000004
$ws.Range("N136").Value = -8184.529500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2139860
$ws.Range("I134").Value = 952.8484999999999
$ws.Range("J134").Value = 5854804.5
$ws.Range("K134").Value = 2858.5455
$ws.Range("L134").Value = 17564413.5
$ws.Range("M134").Value = -323.5454999999997
$ws.Range("N134").Value = -17569483.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1297.3472
$ws.Range("I31").Value = 1007.1316
$ws.Range("J31").Value = 1621.7059
$ws.Range("K31").Value = 1007.1316
$ws.Range("L31").Value = 1621.7059
$ws.Range("M31").Value = -712.1316
$ws.Range("N31").Value = -2211.7059

$ws.Range("H34").Value = 1297.3472
$ws.Range("I34").Value = 1007.1316
$ws.Range("J34").Value = 1621.7059
$ws.Range("K34").Value = 1007.1316
$ws.Range("L34").Value = 1621.7059
$ws.Range("M34").Value = -805.1316
$ws.Range("N34").Value = -2025.7059

$ws.Range("H58").Value = 25001038
$ws.Range("I58").Value = 43479052
$ws.Range("J58").Value = 1371.8823
$ws.Range("K58").Value = 43479052
$ws.Range("L58").Value = 1371.8823
$ws.Range("M58").Value = -43478849
$ws.Range("N58").Value = -1777.8823

$ws.Range("H132").Value = 5209358.5
$ws.Range("I132").Value = 810.9815
$ws.Range("J132").Value = 33335516
$ws.Range("K132").Value = 2432.9445
$ws.Range("L132").Value = 100006548
$ws.Range("M132").Value = 97.05549999999994
$ws.Range("N132").Value = -100011608

$ws.Range("H134").Value = 1092.8148
$ws.Range("I134").Value = 910.45
$ws.Range("J134").Value = 1613.8572
$ws.Range("K134").Value = 2731.35
$ws.Range("L134").Value = 4841.571599999999
$ws.Range("M134").Value = -196.3500000000004
$ws.Range("N134").Value = -9911.571599999999

$ws.Range("H136").Value = 25001038
$ws.Range("I136").Value = 43479052
$ws.Range("J136").Value = 1371.8823
$ws.Range("K136").Value = 130437156
$ws.Range("L136").Value = 4115.6469
$ws.Range("M136").Value = -130434606
$ws.Range("N136").Value = -9215.6469

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 18928.572
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 18928.572
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 56785.716
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -58157.716

$ws.Range("H65").Value = 18928.572
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 18928.572
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 170357.148
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -177221.148

$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 6000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -6630

$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 6000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -8184

$ws.Range("H76").Value = 1250
$ws.Range("I76").Value = 500
$ws.Range("K76").Value = 1500
$ws.Range("M76").Value = -1117

$ws.Range("H79").Value = 1250
$ws.Range("I79").Value = 500
$ws.Range("K79").Value = 1500
$ws.Range("M79").Value = -174

$ws.Range("H81").Value = 20411012
$ws.Range("I81").Value = 142857140
$ws.Range("J81").Value = 3325.8333
$ws.Range("K81").Value = 428571420
$ws.Range("L81").Value = 9977.499899999999
$ws.Range("M81").Value = -428570297
$ws.Range("N81").Value = -12223.4999

$ws.Range("H84").Value = 20411012
$ws.Range("I84").Value = 142857140
$ws.Range("J84").Value = 3325.8333
$ws.Range("K84").Value = 1285714260
$ws.Range("L84").Value = 29932.4997
$ws.Range("M84").Value = -1285708644
$ws.Range("N84").Value = -41164.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5099.204
$ws.Range("I132").Value = 3538.2856
$ws.Range("J132").Value = 9001.5
$ws.Range("K132").Value = 10614.8568
$ws.Range("L132").Value = 27004.5
$ws.Range("M132").Value = -8084.856800000001
$ws.Range("N132").Value = -32064.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12824668
$ws.Range("I132").Value = 17544958
$ws.Range("J132").Value = 12453
$ws.Range("K132").Value = 52634874
$ws.Range("L132").Value = 37359
$ws.Range("M132").Value = -52632344
$ws.Range("N132").Value = -42419

$ws.Range("H136").Value = 21564678
$ws.Range("I136").Value = 3862522.8
$ws.Range("J136").Value = 62500908
$ws.Range("K136").Value = 11587568.4
$ws.Range("L136").Value = 187502724
$ws.Range("M136").Value = -11585018.4
$ws.Range("N136").Value = -187507824

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6253016
$ws.Range("I136").Value = 8067894.5
$ws.Range("J136").Value = 1767.6111
$ws.Range("K136").Value = 24203683.5
$ws.Range("L136").Value = 5302.8333
$ws.Range("M136").Value = -24201133.5
$ws.Range("N136").Value = -10402.8333
